$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1050893
$ws.Range("D2").Value = 74.7
$ws.Range("E2").Value = 785395
$ws.Range("F2").Value = 10.8
$ws.Range("G2").Value = 113896
$ws.Range("H2").Value = 10.6
$ws.Range("I2").Value = 111194
$ws.Range("J2").Value = 3.8
$ws.Range("K2").Value = 40407
$ws.Range("N2").Value = 1050893

$ws.Range("C3").Value = 826821
$ws.Range("D3").Value = 73
$ws.Range("E3").Value = 603451
$ws.Range("F3").Value = 11.9
$ws.Range("G3").Value = 97997
$ws.Range("H3").Value = 12.1
$ws.Range("I3").Value = 99845
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 25528
$ws.Range("N3").Value = 826821

$ws.Range("C4").Value = 143427
$ws.Range("D4").Value = 81.3
$ws.Range("E4").Value = 116675
$ws.Range("F4").Value = 5.1
$ws.Range("G4").Value = 7382
$ws.Range("H4").Value = 5.4
$ws.Range("I4").Value = 7675
$ws.Range("J4").Value = 8.199999999999999
$ws.Range("K4").Value = 11695
$ws.Range("N4").Value = 143427

$ws.Range("C5").Value = 65585
$ws.Range("D5").Value = 79.3
$ws.Range("E5").Value = 51983
$ws.Range("F5").Value = 12.3
$ws.Range("G5").Value = 8099
$ws.Range("H5").Value = 5.3
$ws.Range("I5").Value = 3469
$ws.Range("J5").Value = 3.1
$ws.Range("K5").Value = 2034
$ws.Range("N5").Value = 65585

$ws.Range("C6").Value = 15060
$ws.Range("D6").Value = 88.2
$ws.Range("E6").Value = 13287
$ws.Range("F6").Value = 2.8
$ws.Range("G6").Value = 418
$ws.Range("H6").Value = 1.4
$ws.Range("I6").Value = 205
$ws.Range("J6").Value = 7.6
$ws.Range("K6").Value = 1150
$ws.Range("N6").Value = 15060
